$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Binary Search 2")

# --- Row 5: "Aggressive Cows - Linear Search" entry ---
# Create the hyperlink first (address doubles as the display text so the
# relationship's "display" attribute ends up holding the target URL, same
# as the pre-existing F3/F4 links), then restore the Hyperlink cell style
# (format-only paste from a sibling cell) and finally set the friendly
# GitHub description text that should actually show in the cell.
$ws.Hyperlinks.Add(
    $ws.Range("F5"),
    "https://github.com/ankurnecessary/dsa/blob/main/2_binarySearch/3_aggressive_cows_linear_search.java",
    "",
    "",
    "https://github.com/ankurnecessary/dsa/blob/main/2_binarySearch/3_aggressive_cows_linear_search.java"
) | Out-Null
$ws.Range("F4").Copy() | Out-Null
$ws.Range("F5").PasteSpecial(-4122) | Out-Null
$ws.Range("F5").Value = "dsa/2_binarySearch/3_aggressive_cows_linear_search.java at main · ankurnecessary/dsa · GitHub"

# --- Row 6: "Aggressive Cows - Binary Search" entry ---
$ws.Hyperlinks.Add(
    $ws.Range("F6"),
    "https://github.com/ankurnecessary/dsa/blob/main/2_binarySearch/4_aggressive_cows_binary_search.java",
    "",
    "",
    "https://github.com/ankurnecessary/dsa/blob/main/2_binarySearch/4_aggressive_cows_binary_search.java"
) | Out-Null
$ws.Range("F4").Copy() | Out-Null
$ws.Range("F6").PasteSpecial(-4122) | Out-Null
$ws.Range("F6").Value = "dsa/2_binarySearch/4_aggressive_cows_binary_search.java at main · ankurnecessary/dsa · GitHub"

$excel.CutCopyMode = 0

# Row heights grow to fit the new wrapped hyperlink text (matches rows already using this height).
$ws.Rows.Item(5).RowHeight = 57.6
$ws.Rows.Item(6).RowHeight = 57.6

# Selection ends on F6, as in the saved file.
$ws.Range("F6").Select() | Out-Null
